$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    # Drop the first four data rows (old rows 2-5); remaining data rows
    # 6-20 shift up to become rows 2-16, and the used range / dimension
    # shrinks from A1:C20 to A1:C16 automatically.
    $ws.Range("A2:A5").EntireRow.Delete()

    # Re-index the "Cutoff" column (A) back to a 0-based sequence for the
    # rows that remain (it was 4..18, now it needs to read 0..14).
    for ($r = 2; $r -le 16; $r++) {
        $ws.Cells.Item($r, 1).Value = $r - 2
    }
}
